# UNI-0001 <I> Fix bug that lets user click on the playing board (and the
# board updates!!) even after game result has been decided.
#
# Adds the first real ticket row to the "Main" sheet of the bug tracker.
# Downstream Summary-sheet COUNTIF()s and charts recompute automatically
# from this new row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Main")

$ws.Range("A2").Value = "Bug"
$ws.Range("B2").Value = "UNI-0001"
$ws.Range("C2").Value = "<I> Fix bug that lets user click on the playing board (and the board updates!!)  even after game result has been decided."
$ws.Range("D2").Value = "rhdelaro"
$ws.Range("E2").Value = "rhdelaro"
$ws.Range("F2").Value = "Low"
$ws.Range("G2").Value = "Open"
$ws.Range("H2").Value = "Unresolved"
$ws.Range("I2").Value = "11/15/2013 22:37:42"
$ws.Range("K2").Value = "Clicking anywhere on the board should not result in any update once player has won/lost."

# Row grows to fit the wrapped Summary/Notes text (matches the author's
# manual resize after typing the long strings above).
$ws.Rows("2:2").RowHeight = 42

# Leave the cursor where the author's last edit (Notes, K2) would put it.
$ws.Range("K3").Select() | Out-Null
